$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# URL: matchsource -> matchsync
$ws.Range("B2").Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/nmdp-ethnicity-codes"

# Experimental value: set to literal text "true" (leading apostrophe forces text, not boolean)
$ws.Range("B7").Formula = "'true"

# Date: updated timestamp
$ws.Range("B8").Value = "2024-02-19T18:37:26-06:00"
